$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column to remain text (avoid Excel auto-number coercion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.563.39'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '3.332.89'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '580.02'
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("D6").Value = '175.77'
$ws.Range("E6").Value = '  -3.48%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '3.330.23'
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("E10").Value = '  -0.40%  '
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("E12").Value = '  -2.23%  '
$ws.Range("D13").Value = '0.0000269'
$ws.Range("E13").Value = '  -2.34%  '
$ws.Range("D14").Value = '664.90'
$ws.Range("E14").Value = '  +3.52%  '
$ws.Range("D15").Value = '3.879.52'
$ws.Range("E15").Value = '  +0.70%  '
$ws.Range("D16").Value = '8.41'
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("D17").Value = '67.697.11'
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").Value = '3.334.58'
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("D20").Value = '17.39'
$ws.Range("E20").Value = '  -1.72%  '
$ws.Range("D21").Value = '10.95'
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").Value = '0.890'
$ws.Range("E22").Value = '  -0.94%  '
$ws.Range("D23").Value = '5.47'
$ws.Range("E23").Value = '  +9.03%  '
$ws.Range("D24").Value = '17.07'
$ws.Range("E24").Value = '  -3.57%  '
$ws.Range("D25").Value = '99.31'
$ws.Range("E25").Value = '  +1.23%  '
$ws.Range("D26").Value = '3.85'
$ws.Range("E26").Value = '  -3.48%  '
$ws.Range("E27").Value = '  -4.66%  '
$ws.Range("D28").Value = '9.31'
$ws.Range("E28").Value = '  -3.17%  '
$ws.Range("D29").Value = '33.59'
$ws.Range("E29").Value = '  +2.27%  '
$ws.Range("D30").Value = '7.42'
$ws.Range("E30").Value = '  +11.45%  '
$ws.Range("D31").Value = '8.44'
$ws.Range("E31").Value = '  -1.59%  '
$ws.Range("D32").Value = '575.11'
$ws.Range("E32").Value = '  -4.68%  '
$ws.Range("D33").Value = '10.97'
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").Value = '0.105'
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("D36").Value = '3.697.16'
$ws.Range("E36").Value = '  -5.67%  '
$ws.Range("D37").Value = '56.60'
$ws.Range("E37").Value = '  +1.42%  '
$ws.Range("D38").Value = '3.37'
$ws.Range("E38").Value = '  -6.13%  '
$ws.Range("D39").Value = '34.47'
$ws.Range("E39").Value = '  +4.56%  '
$ws.Range("E40").Value = '  +2.03%  '
$ws.Range("D41").Value = '2.62'
$ws.Range("E41").Value = '  -2.94%  '
$ws.Range("D42").Value = '3.11'
$ws.Range("E42").Value = '  -5.25%  '
$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").Value = '0.334'
$ws.Range("E43").Value = '  -1.43%  '
$ws.Range("B44").Value = 'PEPE'
$ws.Range("C44").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D44").Value = '0.0₃0668'
$ws.Range("E44").Value = '  -3.08%  '
$ws.Range("D45").Value = '3.28'
$ws.Range("E45").Value = '  -2.10%  '
$ws.Range("D46").Value = '0.0405'
$ws.Range("E46").Value = '  -2.55%  '
$ws.Range("E47").Value = '  +1.58%  '
$ws.Range("E48").Value = '  -0.47%  '
$ws.Range("E49").Value = '  -0.11%  '
$ws.Range("E50").Value = '  +0.66%  '
$ws.Range("D51").Value = '128.60'
$ws.Range("E51").Value = '  -1.83%  '
